$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "平測-IT前置準備"
$ws.Range("F3").Value = "OOOOO1"
$ws.Range("F4").Value = "OOOOO2"
$ws.Range("L3").Value = "OOXX1"
$ws.Range("L4").Value = "OOXX2"

$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("L5").Select()
